$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value  = 1.38
$ws.Range("H2").Value  = 5.25
$ws.Range("I2").Value  = 7
$ws.Range("L2").Value  = 7
$ws.Range("N2").Value  = 17
$ws.Range("O2").Value  = 1.2
$ws.Range("P2").Value  = 4.5
$ws.Range("U2").Value  = 1.91
$ws.Range("V2").Value  = 1.91
$ws.Range("W2").Value  = 8
$ws.Range("X2").Value  = 7
$ws.Range("AC2").Value = 17
$ws.Range("AD2").Value = 10
$ws.Range("AE2").Value = 21
$ws.Range("AG2").Value = 251
$ws.Range("AH2").Value = 19
$ws.Range("AN2").Value = 3.4
$ws.Range("AP2").Value = 17
$ws.Range("AQ2").Value = 17
$ws.Range("BB2").Value = 251

# Row 3 updates
$ws.Range("G3").Value  = 2.15
$ws.Range("I3").Value  = 3.75
$ws.Range("J3").Value  = 3
$ws.Range("M3").Value  = 1.13
$ws.Range("N3").Value  = 6
$ws.Range("AN3").Value = 4
$ws.Range("AO3").Value = 13
$ws.Range("AZ3").Value = 81

# Row 4 updates
$ws.Range("G4").Value  = 2.1
$ws.Range("I4").Value  = 3.4
$ws.Range("J4").Value  = 2.75
$ws.Range("L4").Value  = 3.75
$ws.Range("U4").Value  = 1.7
$ws.Range("V4").Value  = 2.05
$ws.Range("Z4").Value  = 19
$ws.Range("AA4").Value = 17
$ws.Range("AC4").Value = 11
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 13
$ws.Range("AF4").Value = 41
$ws.Range("AQ4").Value = 41
$ws.Range("AS4").Value = 151
$ws.Range("AZ4").Value = 51
